$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 500). The commit bumps that date by one day (2023-09-08
# -> 2023-09-09, serials 45177 -> 45178) for every row.
$ws.Range("C2:C500").Value = 45178
